$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 670
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2888

$ws.Range("H20").Value = 2760.5
$ws.Range("I20").Value = 2760.5
$ws.Range("K20").Value = 2760.5
$ws.Range("M20").Value = -2530.5

$ws.Range("H35").Value = 2760.5
$ws.Range("I35").Value = 2760.5
$ws.Range("K35").Value = 2760.5
$ws.Range("M35").Value = -2381.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 22488.334
$ws.Range("J15").Value = 23386
$ws.Range("L15").Value = 23386
$ws.Range("N15").Value = -23840

$ws.Range("H19").Value = 8158.75
$ws.Range("I19").Value = 4324.2856
$ws.Range("K19").Value = 4324.2856
$ws.Range("M19").Value = -4151.2856

$ws.Range("H22").Value = 174.83333
$ws.Range("I22").Value = 174.83333
$ws.Range("K22").Value = 174.83333
$ws.Range("M22").Value = -1.833329999999989

$ws.Range("H134").Value = 7567.0586
$ws.Range("I134").Value = 6883.75
$ws.Range("J134").Value = 18500
$ws.Range("K134").Value = 20651.25
$ws.Range("L134").Value = 55500
$ws.Range("M134").Value = -18116.25
$ws.Range("N134").Value = -60570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1716.3334
$ws.Range("I16").Value = 1716.3334
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1716.3334
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1429.3334
$ws.Range("N16").ClearContents()

$ws.Range("H58").Value = 1699.7059
$ws.Range("J58").Value = 1870
$ws.Range("L58").Value = 1870
$ws.Range("N58").Value = -2276

$ws.Range("H99").Value = 1946.8572
$ws.Range("I99").Value = 1946.8572
$ws.Range("K99").Value = 1946.8572
$ws.Range("M99").Value = -448.8571999999999

$ws.Range("H113").Value = 1716.3334
$ws.Range("I113").Value = 1716.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1716.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 453.6666
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 1946.8572
$ws.Range("I126").Value = 1946.8572
$ws.Range("K126").Value = 5840.571599999999
$ws.Range("M126").Value = -3370.571599999999

$ws.Range("H132").Value = 3283.6365
$ws.Range("I132").Value = 3404.3333
$ws.Range("J132").Value = 2740.5
$ws.Range("K132").Value = 10212.9999
$ws.Range("L132").Value = 8221.5
$ws.Range("M132").Value = -7682.999899999999
$ws.Range("N132").Value = -13281.5

$ws.Range("H134").Value = 4448.409
$ws.Range("I134").Value = 4203.6665
$ws.Range("K134").Value = 12610.9995
$ws.Range("M134").Value = -10075.9995

$ws.Range("H136").Value = 1699.7059
$ws.Range("J136").Value = 1870
$ws.Range("L136").Value = 5610
$ws.Range("N136").Value = -10710

$ws.Range("H112").Value = 2241
$ws.Range("I112").Value = 4427
$ws.Range("J112").Value = 55
$ws.Range("K112").Value = 13281
$ws.Range("L112").Value = 165
$ws.Range("M112").Value = -12173
$ws.Range("N112").Value = -2381

$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 814.6
$ws.Range("J34").Value = 1198
$ws.Range("L34").Value = 3594
$ws.Range("N34").Value = -3762

$ws.Range("H37").Value = 99999.4
$ws.Range("J37").Value = 99999.4
$ws.Range("L37").Value = 299998.2
$ws.Range("N37").Value = -300222.2

$ws.Range("H50").Value = 875.6667
$ws.Range("I50").Value = 326.25
$ws.Range("K50").Value = 978.75
$ws.Range("M50").Value = -497.75

$ws.Range("H53").Value = 875.6667
$ws.Range("I53").Value = 326.25
$ws.Range("K53").Value = 978.75
$ws.Range("M53").Value = -497.75

$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("K87").Value = 3000
$ws.Range("M87").Value = -1752

$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("K90").Value = 9000
$ws.Range("M90").Value = -2760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 525.1875
$ws.Range("I2").Value = 487.8
$ws.Range("K2").Value = 487.8
$ws.Range("M2").Value = -374.8

$ws.Range("H15").Value = 7700
$ws.Range("J15").Value = 7700
$ws.Range("L15").Value = 7700
$ws.Range("N15").Value = -8276

$ws.Range("H81").Value = 7700
$ws.Range("J81").Value = 7700
$ws.Range("L81").Value = 7700
$ws.Range("N81").Value = -9696

$ws.Range("H84").Value = 7700
$ws.Range("J84").Value = 7700
$ws.Range("L84").Value = 23100
$ws.Range("N84").Value = -33084

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2285.8572
$ws.Range("I22").Value = 2500.1667
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2500.1667
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -2205.1667
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 2285.8572
$ws.Range("I27").Value = 2500.1667
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2500.1667
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2393.1667
$ws.Range("N27").Value = -1214

$ws.Range("H136").Value = 3353.2173
$ws.Range("I136").Value = 2971.3076
$ws.Range("K136").Value = 8913.9228
$ws.Range("M136").Value = -6363.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4572.1665
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4572.1665
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988

$ws.Range("H126").Value = 38922.8
$ws.Range("I126").Value = 34560.215
$ws.Range("J126").Value = 99999
$ws.Range("K126").Value = 103680.645
$ws.Range("L126").Value = 299997
$ws.Range("M126").Value = -101210.645
$ws.Range("N126").Value = -304937

$ws.Range("H136").Value = 5607.579
$ws.Range("J136").Value = 5431.25
$ws.Range("L136").Value = 16293.75
$ws.Range("N136").Value = -21393.75
